{"js": "// Update the worksheet date and all twenty-six \"dividend\u00f7divisor=\" problems.\n// Every old value below is unique within the document, so a simple\n// search-and-replace-first-match per pair is unambiguous and safe to run\n// in sequence (later pairs never collide with not-yet-processed earlier\n// pairs because all values, old and new, are distinct strings).\nconst replacements = [\n  [\"2024-02-04 Sunday\", \"2024-02-05 Monday\"],\n  [\"124\u00f74=\", \"260\u00f76=\"],\n  [\"740\u00f78=\", \"120\u00f75=\"],\n  [\"876\u00f75=\", \"433\u00f72=\"],\n  [\"501\u00f75=\", \"462\u00f73=\"],\n  [\"406\u00f74=\", \"202\u00f73=\"],\n  [\"764\u00f75=\", \"878\u00f74=\"],\n  [\"259\u00f77=\", \"353\u00f79=\"],\n  [\"525\u00f77=\", \"286\u00f72=\"],\n  [\"202\u00f76=\", \"320\u00f72=\"],\n  [\"754\u00f73=\", \"179\u00f75=\"],\n  [\"420\u00f76=\", \"197\u00f78=\"],\n  [\"151\u00f73=\", \"235\u00f73=\"],\n  [\"960\u00f72=\", \"783\u00f78=\"],\n  [\"565\u00f77=\", \"372\u00f78=\"],\n  [\"188\u00f79=\", \"994\u00f77=\"],\n  [\"532\u00f79=\", \"374\u00f76=\"],\n  [\"208\u00f74=\", \"485\u00f75=\"],\n  [\"402\u00f78=\", \"947\u00f76=\"],\n  [\"908\u00f72=\", \"869\u00f78=\"],\n  [\"965\u00f72=\", \"613\u00f75=\"],\n  [\"550\u00f73=\", \"705\u00f77=\"],\n  [\"946\u00f77=\", \"651\u00f77=\"],\n  [\"118\u00f79=\", \"475\u00f72=\"],\n  [\"637\u00f78=\", \"162\u00f77=\"],\n  [\"738\u00f77=\", \"402\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Update the worksheet date and all twenty-six \"dividend\u00f7divisor=\" problems.\n# Every \"old\" value is unique within the document, so Find/Replace on the\n# whole-document range for each pair, in order, is unambiguous. (The final\n# pair's replacement text, \"402\u00f78=\", happens to equal an earlier pair's old\n# text, but that earlier occurrence was already consumed by the time we get\n# here, so no collision occurs.)\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-04 Sunday\", \"2024-02-05 Monday\"),\n    @(\"124\u00f74=\", \"260\u00f76=\"),\n    @(\"740\u00f78=\", \"120\u00f75=\"),\n    @(\"876\u00f75=\", \"433\u00f72=\"),\n    @(\"501\u00f75=\", \"462\u00f73=\"),\n    @(\"406\u00f74=\", \"202\u00f73=\"),\n    @(\"764\u00f75=\", \"878\u00f74=\"),\n    @(\"259\u00f77=\", \"353\u00f79=\"),\n    @(\"525\u00f77=\", \"286\u00f72=\"),\n    @(\"202\u00f76=\", \"320\u00f72=\"),\n    @(\"754\u00f73=\", \"179\u00f75=\"),\n    @(\"420\u00f76=\", \"197\u00f78=\"),\n    @(\"151\u00f73=\", \"235\u00f73=\"),\n    @(\"960\u00f72=\", \"783\u00f78=\"),\n    @(\"565\u00f77=\", \"372\u00f78=\"),\n    @(\"188\u00f79=\", \"994\u00f77=\"),\n    @(\"532\u00f79=\", \"374\u00f76=\"),\n    @(\"208\u00f74=\", \"485\u00f75=\"),\n    @(\"402\u00f78=\", \"947\u00f76=\"),\n    @(\"908\u00f72=\", \"869\u00f78=\"),\n    @(\"965\u00f72=\", \"613\u00f75=\"),\n    @(\"550\u00f73=\", \"705\u00f77=\"),\n    @(\"946\u00f77=\", \"651\u00f77=\"),\n    @(\"118\u00f79=\", \"475\u00f72=\"),\n    @(\"637\u00f78=\", \"162\u00f77=\"),\n    @(\"738\u00f77=\", \"402\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
